# User CRUD complete + Email portal fixes
# Populates rows 4-23 with the same sender/snippet/subject/fake pattern
# already present in rows 2-3, wires up a mailto hyperlink on each new
# "sender" cell (column A), and updates the sheet's view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two alternating "template" rows already on the sheet.
$senders  = @("me@me.me", "you@you.you")
$snippets = @("hello", "goodbye")
$subjects = @("Important", "Not")
$fakes    = @($false, $true)

for ($row = 4; $row -le 23; $row++) {
    $idx = ($row - 4) % 2

    $ws.Cells.Item($row, 1).Value = $senders[$idx]
    $ws.Cells.Item($row, 2).Value = $snippets[$idx]
    $ws.Cells.Item($row, 3).Value = $subjects[$idx]
    $ws.Cells.Item($row, 4).Value = $fakes[$idx]

    $anchor = $ws.Cells.Item($row, 1)
    $ws.Hyperlinks.Add($anchor, "mailto:" + $senders[$idx]) | Out-Null
}

# Hyperlinks.Add() stamps a fresh style index on each cell it touches;
# re-apply the workbook's built-in "Hyperlink" cell style in one shot so
# the new cells line up with the existing A2/A3 styling.
$ws.Range("A4:A23").Style = "Hyperlink"

# Match the saved view: scrolled down with F17 selected.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F17").Select() | Out-Null
